$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data through row 115. Append 14 more rows
# (116-129) that repeat the same two-row "BBQ" transaction pattern
# already present in rows 114-115, seven times in a row.
# Using Copy/Paste (instead of typing literal values) keeps the
# numeric-looking text cells (column C) stored as text/shared-strings,
# exactly like the source rows, instead of being auto-converted to
# numbers.

$source = $ws.Range("A114:G115")

for ($i = 0; $i -lt 7; $i++) {
    $destStartRow = 116 + ($i * 2)
    $destEndRow = $destStartRow + 1
    $dest = $ws.Range("A" + $destStartRow + ":G" + $destEndRow)
    $source.Copy($dest)
}

$wb.Save()
